$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.469106674194336
$ws.Range("B1").Value = 3.989879369735718
$ws.Range("C1").Value = 3.478955268859863
$ws.Range("D1").Value = 1.914044976234436
$ws.Range("E1").Value = 0.9573580026626587
